# Update the "cryptos" worksheet (coinranking.com scrape) with the latest
# price / 1h-volume figures, as produced by the Thu Aug 24 09:36:33 UTC 2023
# GitHub Actions run.
#
# Column D ("Price") and column E ("Volume(1h)") hold plain text in this
# workbook (not numbers/percentages) -- several prices parse as valid
# numbers (e.g. "1.002", "0.5300"), so a naive `.Value = ...` assignment
# would have Excel silently reinterpret them as numeric and, worse, drop
# meaningful trailing zeros (e.g. "0.5300" -> 0.53). To avoid that we
# briefly force the cell to Text format before writing, then restore its
# original ("Normal") style so formatting is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2 = @("26.523.60", "  +1.80%  ")
    3 = @("1.674.28", "  +1.74%  ")
    4 = @("1.002", "  +0.02%  ")
    5 = @("220.16", "  +2.21%  ")
    6 = @("0.5300", "  +1.45%  ")
    7 = @($null, "  +0.04%  ")
    8 = @("0.2682", "  +2.65%  ")
    9 = @($null, "  +0.42%  ")
    10 = @("21.85", "  +4.79%  ")
    11 = @("0.07804", "  +1.86%  ")
    12 = @("1.677.08", "  +1.79%  ")
    13 = @("4.500", "  +1.73%  ")
    14 = @("0.5581", "  +0.51%  ")
    15 = @("0.0₅8337", "  +0.70%  ")
    16 = @($null, "  +1.26%  ")
    17 = @("26.544.43", "  +1.82%  ")
    18 = @("1.002", "  +0.04%  ")
    19 = @("4.777", "  +1.12%  ")
    20 = @("193.55", "  +2.75%  ")
    21 = @($null, "  +1.53%  ")
    22 = @("6.316", "  +1.10%  ")
    23 = @("1.003", "  +0.06%  ")
    24 = @("0.1272", "  +4.31%  ")
    25 = @("138.69", $null)
    26 = @("7.405", "  -0.16%  ")
    27 = @("16.33", "  +3.07%  ")
    28 = @("1.429", "  +3.15%  ")
    29 = @("0.06264", "  +5.03%  ")
    30 = @("1.289", "  +1.97%  ")
    31 = @($null, "  +6.13%  ")
    32 = @("3.420", "  +0.62%  ")
    33 = @("1.695", "  +2.38%  ")
    34 = @("1.011", "  +1.33%  ")
    35 = @("0.6198", "  +10.22%  ")
    36 = @($null, "  +1.21%  ")
    37 = @("2.789", "  +1.28%  ")
    38 = @("0.01619", "  +0.63%  ")
    39 = @("6.056", "  +3.72%  ")
    40 = @("1.095.29", "  +6.63%  ")
    41 = @("0.8606", "  +0.43%  ")
    42 = @($null, "  +0.00%  ")
    43 = @("100.64", "  +1.46%  ")
    44 = @("1.821.06", "  +1.42%  ")
    45 = @("59.01", "  +5.78%  ")
    46 = @("0.0₈109", "  -1.88%  ")
    47 = @("8.200", "  +1.88%  ")
    48 = @("1.530", "  +10.78%  ")
    49 = @("0.9969", "  -0.71%  ")
    50 = @("0.05198", "  +0.95%  ")
    51 = @("6.021", "  +1.43%  ")
}

foreach ($row in $changes.Keys) {
    $newPrice = $changes[$row][0]
    $newVolume = $changes[$row][1]

    if ($null -ne $newPrice) {
        $priceCell = $ws.Cells.Item($row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $newPrice
        $priceCell.Style = "Normal"
    }

    if ($null -ne $newVolume) {
        $ws.Cells.Item($row, 5).Value = $newVolume
    }
}
